$wb = $excel.ActiveWorkbook

function Add-FanSheet {
    param(
        [string]$Name,
        [double]$Pwm100,
        [double]$Sound100,
        [double]$PwmMin,
        [double]$SoundMin
    )

    $after = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $after)
    $ws.Name = $Name

    $ws.Range("A1").Value = "Model"
    $ws.Range("B1").Value = "PWM - [%]"
    $ws.Range("C1").Value = "Sound Level - [db-A]"

    $ws.Rows.Item(1).RowHeight = 43.5
    $ws.Range("C1").HorizontalAlignment = -4108
    $ws.Range("C1").VerticalAlignment = -4108
    $ws.Range("C1").WrapText = $true

    $ws.Range("A2").Value = $Name
    $ws.Range("B2").Value = $Pwm100
    $ws.Range("C2").Value = $Sound100

    $ws.Range("A3").Value = $Name
    $ws.Range("B3").Value = $PwmMin
    $ws.Range("C3").Value = $SoundMin

    $ws.Range("A1").Select()

    return $ws
}

$sheet1 = Add-FanSheet "9GT1224P1S001" 100 58 35 41
$sheet2 = Add-FanSheet "9GT0924P1M001" 100 58 30 37

$sheet2.Range("C4").Select()

Write-Output "Added sheets 9GT1224P1S001 and 9GT0924P1M001"
